$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Amh"
$ws.Cells.Item(2,3).Value = "Acvr1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.363728
$ws.Cells.Item(2,8).Value = 1.091184
$ws.Cells.Item(2,9).Value = 0.3031196334509865
$ws.Cells.Item(2,10).Value = 0.3031196334509865
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 4.695610666666666
$ws.Cells.Item(2,14).Value = 14.086832
$ws.Cells.Item(2,15).Value = 0.1802066564018305
$ws.Cells.Item(2,16).Value = 0.1802066564018305
$ws.Cells.Item(2,17).Value = 1.707925076565333
$ws.Cells.Item(2,18).Value = 15.371325689088
$ws.Cells.Item(2,19).Value = 0.05462417563395074
$ws.Cells.Item(2,20).Value = 0.05462417563395073

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Amh"
$ws.Cells.Item(3,3).Value = "Acvr1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.363728
$ws.Cells.Item(3,8).Value = 1.091184
$ws.Cells.Item(3,9).Value = 0.3031196334509865
$ws.Cells.Item(3,10).Value = 0.3031196334509865
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 15.51448033333333
$ws.Cells.Item(3,14).Value = 46.543441
$ws.Cells.Item(3,15).Value = 0.5954098039960916
$ws.Cells.Item(3,16).Value = 0.5954098039960916
$ws.Cells.Item(3,17).Value = 5.643050902682667
$ws.Cells.Item(3,18).Value = 50.787458124144
$ws.Cells.Item(3,19).Value = 0.180480401540419
$ws.Cells.Item(3,20).Value = 0.180480401540419

# Row 4: ECs -> sCs
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Amh"
$ws.Cells.Item(4,3).Value = "Acvr1"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.363728
$ws.Cells.Item(4,8).Value = 1.091184
$ws.Cells.Item(4,9).Value = 0.3031196334509865
$ws.Cells.Item(4,10).Value = 0.3031196334509865
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 5.846719333333333
$ws.Cells.Item(4,14).Value = 17.540158
$ws.Cells.Item(4,15).Value = 0.2243835396020779
$ws.Cells.Item(4,16).Value = 0.2243835396020779
$ws.Cells.Item(4,17).Value = 2.126615529674666
$ws.Cells.Item(4,18).Value = 19.139539767072
$ws.Cells.Item(4,19).Value = 0.06801505627661678
$ws.Cells.Item(4,20).Value = 0.06801505627661676

# Row 5: FAPs -> ECs
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Amh"
$ws.Cells.Item(5,3).Value = "Acvr1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.520663
$ws.Cells.Item(5,8).Value = 1.561989
$ws.Cells.Item(5,9).Value = 0.4339043947991109
$ws.Cells.Item(5,10).Value = 0.4339043947991109
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 4.695610666666666
$ws.Cells.Item(5,14).Value = 14.086832
$ws.Cells.Item(5,15).Value = 0.1802066564018305
$ws.Cells.Item(5,16).Value = 0.1802066564018305
$ws.Cells.Item(5,17).Value = 2.444830736538667
$ws.Cells.Item(5,18).Value = 22.003476628848
$ws.Cells.Item(5,19).Value = 0.07819246018480758
$ws.Cells.Item(5,20).Value = 0.07819246018480756

# Row 6: FAPs -> FAPs
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Amh"
$ws.Cells.Item(6,3).Value = "Acvr1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.520663
$ws.Cells.Item(6,8).Value = 1.561989
$ws.Cells.Item(6,9).Value = 0.4339043947991109
$ws.Cells.Item(6,10).Value = 0.4339043947991109
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 15.51448033333333
$ws.Cells.Item(6,14).Value = 46.543441
$ws.Cells.Item(6,15).Value = 0.5954098039960916
$ws.Cells.Item(6,16).Value = 0.5954098039960916
$ws.Cells.Item(6,17).Value = 8.077815873794332
$ws.Cells.Item(6,18).Value = 72.70034286414901
$ws.Cells.Item(6,19).Value = 0.2583509306603813
$ws.Cells.Item(6,20).Value = 0.2583509306603813

# Row 7: FAPs -> sCs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Amh"
$ws.Cells.Item(7,3).Value = "Acvr1"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.520663
$ws.Cells.Item(7,8).Value = 1.561989
$ws.Cells.Item(7,9).Value = 0.4339043947991109
$ws.Cells.Item(7,10).Value = 0.4339043947991109
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 5.846719333333333
$ws.Cells.Item(7,14).Value = 17.540158
$ws.Cells.Item(7,15).Value = 0.2243835396020779
$ws.Cells.Item(7,16).Value = 0.2243835396020779
$ws.Cells.Item(7,17).Value = 3.044170428251333
$ws.Cells.Item(7,18).Value = 27.397533854262
$ws.Cells.Item(7,19).Value = 0.09736100395392194
$ws.Cells.Item(7,20).Value = 0.09736100395392193

# Row 8: sCs -> ECs
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Amh"
$ws.Cells.Item(8,3).Value = "Acvr1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.3155576666666667
$ws.Cells.Item(8,8).Value = 0.9466730000000001
$ws.Cells.Item(8,9).Value = 0.2629759717499027
$ws.Cells.Item(8,10).Value = 0.2629759717499027
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 4.695610666666666
$ws.Cells.Item(8,14).Value = 14.086832
$ws.Cells.Item(8,15).Value = 0.1802066564018305
$ws.Cells.Item(8,16).Value = 0.1802066564018305
$ws.Cells.Item(8,17).Value = 1.481735945548444
$ws.Cells.Item(8,18).Value = 13.335623509936
$ws.Cells.Item(8,19).Value = 0.0473900205830722
$ws.Cells.Item(8,20).Value = 0.04739002058307219

# Row 9: sCs -> FAPs
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Amh"
$ws.Cells.Item(9,3).Value = "Acvr1"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.3155576666666667
$ws.Cells.Item(9,8).Value = 0.9466730000000001
$ws.Cells.Item(9,9).Value = 0.2629759717499027
$ws.Cells.Item(9,10).Value = 0.2629759717499027
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 15.51448033333333
$ws.Cells.Item(9,14).Value = 46.543441
$ws.Cells.Item(9,15).Value = 0.5954098039960916
$ws.Cells.Item(9,16).Value = 0.5954098039960916
$ws.Cells.Item(9,17).Value = 4.895713213532556
$ws.Cells.Item(9,18).Value = 44.06141892179301
$ws.Cells.Item(9,19).Value = 0.1565784717952912
$ws.Cells.Item(9,20).Value = 0.1565784717952912

# Row 10: sCs -> sCs
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Amh"
$ws.Cells.Item(10,3).Value = "Acvr1"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.3155576666666667
$ws.Cells.Item(10,8).Value = 0.9466730000000001
$ws.Cells.Item(10,9).Value = 0.2629759717499027
$ws.Cells.Item(10,10).Value = 0.2629759717499027
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 5.846719333333333
$ws.Cells.Item(10,14).Value = 17.540158
$ws.Cells.Item(10,15).Value = 0.2243835396020779
$ws.Cells.Item(10,16).Value = 0.2243835396020779
$ws.Cells.Item(10,17).Value = 1.844977110481556
$ws.Cells.Item(10,18).Value = 16.604793994334
$ws.Cells.Item(10,19).Value = 0.0590074793715392
$ws.Cells.Item(10,20).Value = 0.0590074793715392
